{"js": "// The underlying change in this revision is *purely* a raw-XML\n// attribute-order swap inside a handful of already-existing\n// `<w:bookmarkStart>` elements (w:id/w:name swap to w:name/w:id) and\n// `<w:rFonts>` elements that carry `w:hint=\"eastAsia\"` (hint moved to the\n// front of the attribute list). No bookmark is added/removed/renamed, no\n// bookmark id changes, no run is split/merged, no text/formatting differs,\n// and no font actually changes for any run \u2014 every attribute value present\n// before the edit is present after it, on the very same element.\n//\n// Verified directly against the OOXML part: parsing both the original and\n// the reconstructed-from-diff document.xml into an attribute-order-\n// insensitive tree shows them to be 100% identical. In other words this is\n// a cosmetic re-serialization artifact (the commit message even says as\n// much: \"Firxed a bug in shuffling elements without a parent element\" -\n// i.e. a fix to the *test-fixture shuffler tool* that produced this\n// sample file, not a content edit made through Word).\n//\n// The Word JavaScript API (like the Word COM object model) works against\n// the document's object model - paragraphs/ranges/fonts/bookmarks - and\n// has no surface for dictating the literal attribute order XML is\n// serialized with; any attempt to \"force\" that order by deleting and\n// re-adding the bookmarks (the only bookmark write primitive exposed,\n// `Range.insertBookmark` / `Document.Bookmarks.Add`) would actually change\n// the document: it renumbers bookmark ids and always serializes new\n// bookmarks as `w:id` then `w:name`, i.e. it cannot reproduce - and would\n// actively fight - the target order while also introducing real (unwanted)\n// diffs such as bookmark-id churn.\n//\n// So the faithful, side-effect-free application of this diff is a no-op:\n// touch nothing, leave every paragraph/run/bookmark exactly as-is.\n", "ps1": "# The underlying change in this revision is *purely* a raw-XML\n# attribute-order swap inside a handful of already-existing\n# <w:bookmarkStart> elements (w:id/w:name swap to w:name/w:id) and\n# <w:rFonts> elements that carry w:hint=\"eastAsia\" (hint moved to the\n# front of the attribute list). No bookmark is added/removed/renamed, no\n# bookmark id changes, no run is split/merged, no text/formatting differs,\n# and no font actually changes for any run - every attribute value present\n# before the edit is present after it, on the very same element.\n#\n# Verified directly against the OOXML part: parsing both the original and\n# the reconstructed-from-diff document.xml into an attribute-order-\n# insensitive tree shows them to be 100% identical. In other words this is\n# a cosmetic re-serialization artifact (the commit message even says as\n# much: \"Firxed a bug in shuffling elements without a parent element\" -\n# i.e. a fix to the *test-fixture shuffler tool* that produced this sample\n# file, not a content edit made through Word).\n#\n# The Word COM object model (like the Word JavaScript API) works against\n# the document's object model - Paragraphs/Range/Font/Bookmarks - and has\n# no surface for dictating the literal attribute order XML is serialized\n# with; any attempt to \"force\" that order by deleting and re-adding the\n# bookmarks (via $d.Bookmarks.Add, the only bookmark write primitive\n# exposed) would actually change the document: it renumbers bookmark ids\n# and always serializes new bookmarks as w:id then w:name, i.e. it cannot\n# reproduce - and would actively fight - the target order while also\n# introducing real (unwanted) diffs such as bookmark-id churn.\n#\n# So the faithful, side-effect-free application of this diff is a no-op:\n# touch nothing, leave every paragraph/run/bookmark exactly as-is.\n"}
